$d = $word.ActiveDocument

# --- Locate anchor paragraphs by their text ---
# NB: Range.Text includes the trailing paragraph-mark character (CR),
# so trim it before comparing.
function Find-ParaIndex($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $t = $doc.Paragraphs.Item($i).Range.Text.TrimEnd([char]13)
        if ($t -eq $text) { return $i }
    }
    return -1
}

$assumptionsIdx = Find-ParaIndex $d "Assumptions:"

# ---------------------------------------------------------------------
# 1) Insert the three new "Assumptions" bullet paragraphs right after
#    the "Assumptions:" heading, reusing the first blank paragraph that
#    currently follows it as the insertion point (so we don't add an
#    extra paragraph mark).
# ---------------------------------------------------------------------
$firstBlank = $d.Paragraphs.Item($assumptionsIdx + 1)
$insertionPoint = $firstBlank.Range
$insertionPoint.Collapse(1)

$newParasXml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
    '<w:p><w:r><w:t>All hardware works as advertised and without any issues</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>' + `
    '<w:p><w:r><w:t>No change of staff during the entirety of the project</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>' + `
    '<w:p><w:r><w:t>Scope of work can change at any point.</w:t></w:r></w:p>' + `
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertionPoint.InsertXML($newParasXml) | Out-Null

# ---------------------------------------------------------------------
# 2) Only a single blank paragraph should remain between the new bullet
#    list and the "Cost:" heading (there used to be three). Delete the
#    extra blank paragraphs, keeping just one.
# ---------------------------------------------------------------------
$costIdx = Find-ParaIndex $d "Cost:"

while (($costIdx - 1) -ge 1 -and `
       $d.Paragraphs.Item($costIdx - 1).Range.Text.TrimEnd([char]13) -eq "" -and `
       $d.Paragraphs.Item($costIdx - 2).Range.Text.TrimEnd([char]13) -eq "") {
    $extraBlank = $d.Paragraphs.Item($costIdx - 1)
    $extraBlank.Range.Delete()
    $costIdx = $costIdx - 1
}

# ---------------------------------------------------------------------
# 3) Rework the "Cost:" paragraph and the cost-detail paragraph that
#    follows it: the bookmark moves to the end of the "Cost:" paragraph
#    and the leading spaces + sentence become a single run.
# ---------------------------------------------------------------------
$costPara = $d.Paragraphs.Item($costIdx)
$detailPara = $d.Paragraphs.Item($costIdx + 1)
$combined = $d.Range($costPara.Range.Start, $detailPara.Range.End)

$costXml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
    '<w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Cost:</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>' + `
    '<w:p><w:r><w:t xml:space="preserve">   The majority of cost for this project will be labor cost.</w:t></w:r></w:p>' + `
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$combined.InsertXML($costXml) | Out-Null

Write-Output "done"
